$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O7").Value = 0.01567292213439941
$ws.Range("O8").Value = 0.001512765884399414
$ws.Range("O10").Value = 0
$ws.Range("O11").Value = 0.02133274078369141
$ws.Range("O12").Value = 0.02607297897338867
$ws.Range("O13").Value = 0.01903104782104492
$ws.Range("O14").Value = 0.3071813583374023
$ws.Range("O15").Value = 0
$ws.Range("O16").Value = 0.003284215927124023
$ws.Range("O17").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("O20").Value = 0.001423835754394531
$ws.Range("O21").Value = 0.002333164215087891
$ws.Range("O22").Value = 0.01639580726623535
$ws.Range("O23").Value = 0
$ws.Range("O25").Value = 0.04963326454162598
$ws.Range("O26").Value = 2.335688352584839
$ws.Range("O28").Value = 0.01984119415283203
$ws.Range("O30").Value = 0.001100778579711914
$ws.Range("O31").Value = 0.006756782531738281
$ws.Range("O32").Value = 0.0005242824554443359
$ws.Range("O33").Value = 0.01284980773925781
$ws.Range("O34").Value = 0.1166512966156006
$ws.Range("O35").Value = 0.06656217575073242
$ws.Range("O37").Value = 0.124117374420166
$ws.Range("O38").Value = 0.07082605361938477
$ws.Range("O39").Value = 0.05001258850097656
$ws.Range("O40").Value = 0.001550436019897461
$ws.Range("O41").Value = 0.04971671104431152
$ws.Range("O42").Value = 0.01993155479431152
$ws.Range("O43").Value = 0.06643915176391602
$ws.Range("O44").Value = 0
$ws.Range("O46").Value = 0.05138611793518066
$ws.Range("O47").Value = 0.04651403427124023
$ws.Range("O48").Value = 0.0432884693145752
$ws.Range("O49").Value = 0.04978203773498535
$ws.Range("O51").Value = 0.01639413833618164
$ws.Range("O52").Value = 0.001251220703125
$ws.Range("O53").Value = 0.01732349395751953
$ws.Range("O54").Value = 0.01561999320983887
$ws.Range("O55").Value = 0
$ws.Range("O57").Value = 0.003082990646362305
$ws.Range("O58").Value = 0.03150129318237305
$ws.Range("O59").Value = 0.005504131317138672
$ws.Range("O60").Value = 0
$ws.Range("O62").Value = 0.008042573928833008
$ws.Range("O63").Value = 0.007013082504272461
$ws.Range("O64").Value = 0.01845765113830566
$ws.Range("O65").Value = 0.7435266971588135
$ws.Range("O66").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("O71").Value = 0.5749311447143555
$ws.Range("O72").Value = 6.437307834625244
$ws.Range("O73").Value = 0.02820682525634766
$ws.Range("O78").Value = 0.00199437141418457
$ws.Range("O79").Value = 0.009083986282348633
$ws.Range("O80").Value = 0.009607553482055664
$ws.Range("O81").Value = 0.1135287284851074
$ws.Range("O82").Value = 0.007151126861572266
$ws.Range("O83").Value = 0.003013372421264648
$ws.Range("O84").Value = 4.214487314224243
$ws.Range("O85").Value = 0.001000165939331055
$ws.Range("O86").Value = 0.006630659103393555
$ws.Range("O87").Value = 0
$ws.Range("O88").Value = 0.003081083297729492
$ws.Range("O89").Value = 0.04815411567687988
$ws.Range("O90").Value = 8.103768348693848
$ws.Range("O91").Value = 0.06442689895629883
$ws.Range("O92").Value = 0.07472062110900879
$ws.Range("O93").Value = 0.0462651252746582
$ws.Range("O94").Value = 0.002055168151855469
$ws.Range("O95").Value = 0.001998186111450195
$ws.Range("O96").Value = 0.07336950302124023
$ws.Range("O97").Value = 0.006829023361206055
$ws.Range("O98").Value = 0.007189750671386719
$ws.Range("O99").Value = 0.05993318557739258
$ws.Range("O100").Value = 0.08552718162536621
$ws.Range("O101").Value = 0.007103919982910156
$ws.Range("O103").Value = 0.03928112983703613
$ws.Range("O104").Value = 0.07538104057312012
$ws.Range("O105").Value = 0.002104759216308594
$ws.Range("O106").Value = 0.01195645332336426
$ws.Range("O107").Value = 0.01785564422607422
$ws.Range("O108").Value = 0.07224702835083008
$ws.Range("O109").Value = 0.02356219291687012
$ws.Range("O110").Value = 0.05336976051330566
$ws.Range("O111").Value = 0.004682064056396484
$ws.Range("O112").Value = 0
$ws.Range("O114").Value = 0.004708051681518555
$ws.Range("O115").Value = 0.03999662399291992
$ws.Range("O116").Value = 0.06065940856933594
$ws.Range("O117").Value = 0.001003265380859375
$ws.Range("O119").Value = 0.001003026962280273
$ws.Range("O122").Value = 0.0009982585906982422
$ws.Range("O124").Value = 0.001022100448608398
$ws.Range("O127").Value = 0.0009698867797851562
$ws.Range("O128").Value = 0.00108790397644043
$ws.Range("O133").Value = 0.00099945068359375
$ws.Range("O134").Value = 0.00153803825378418
$ws.Range("O135").Value = 0.002008438110351562
$ws.Range("O136").Value = 0.001044511795043945
$ws.Range("O138").Value = 0.005125522613525391
$ws.Range("O139").Value = 0.001000165939331055
$ws.Range("O140").Value = 0.01351141929626465
$ws.Range("O141").Value = 0.002999544143676758
$ws.Range("O142").Value = 0
$ws.Range("O143").Value = 0.001908063888549805
$ws.Range("O146").Value = 0.003236770629882812
$ws.Range("O147").Value = 0.0374910831451416
$ws.Range("O148").Value = 0.03570032119750977
$ws.Range("O149").Value = 0.1536204814910889
$ws.Range("O150").Value = 4.650670528411865
$ws.Range("O151").Value = 0.35862135887146
$ws.Range("O153").Value = 0
$ws.Range("O156").Value = 0.0009999275207519531
$ws.Range("O157").Value = 0.003350496292114258
$ws.Range("O159").Value = 0.006999492645263672
$ws.Range("O160").Value = 0.001908063888549805
$ws.Range("O161").Value = 0.0008993148803710938
$ws.Range("O162").Value = 0.00293278694152832
$ws.Range("O163").Value = 0
$ws.Range("O164").Value = 0.03308963775634766
$ws.Range("O165").Value = 0.002029657363891602
$ws.Range("O166").Value = 0.01582670211791992
$ws.Range("O168").Value = 0.02868056297302246
$ws.Range("O169").Value = 0.0169224739074707
$ws.Range("O170").Value = 0.02729916572570801
$ws.Range("O171").Value = 0.01534414291381836
$ws.Range("O172").Value = 0.002001047134399414
$ws.Range("O173").Value = 0.01712250709533691
$ws.Range("O174").Value = 0.008299350738525391
$ws.Range("O176").Value = 0
$ws.Range("O177").Value = 0.00311732292175293
